# 7.5.10 - bug fixes
#
# The "advanced_tags" sheet's sample/demo rows are cleared out (keeping just
# the header row and the formatting on the now-empty rows below it), the
# "advanced_tags" tab becomes the active/selected tab (instead of
# "tags_fqdn"), its view is zoomed to 200%, and its selection is set to the
# old data range.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("advanced_tags")

# Clear out all the sample rows (A2:F18) but leave the row formatting (style)
# that already exists on A7:A18 in place - this also shrinks the shared
# string table down to just what's still referenced by the workbook.
$ws.Range("A2:F18").ClearContents()

# Make "advanced_tags" the active sheet/tab (was "tags_fqdn").
$ws.Activate()

# Restore the selection to cover the old data block.
[void]$ws.Range("A2:F18").Select()

# Zoom the view to 200%.
$excel.ActiveWindow.Zoom = 200
